$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number need to be forced to Text
# so Excel does not auto-convert the inline string into a numeric cell,
# matching the source workbook which stores every D/E cell as text.
$textCells = @(
    "D5",
    "D6",
    "D10",
    "D13",
    "D17",
    "D19",
    "D20",
    "D22",
    "D23",
    "D28",
    "D30",
    "D31",
    "D32",
    "D33",
    "D35",
    "D37",
    "D39",
    "D40",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.060.50'
$ws.Range("D3").Value = '2.312.52'
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '301.89'
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").Value = '98.74'
$ws.Range("E6").Value = '  -2.23%  '
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  +2.21%  '
$ws.Range("D10").Value = '35.71'
$ws.Range("E10").Value = '  +1.92%  '
$ws.Range("E11").Value = '  -1.08%  '
$ws.Range("E12").Value = '  -1.06%  '
$ws.Range("D13").Value = '17.82'
$ws.Range("E13").Value = '  -1.60%  '
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("D15").Value = '2.670.89'
$ws.Range("E15").Value = '  -0.69%  '
$ws.Range("D16").Value = '2.330.93'
$ws.Range("E16").Value = '  +2.78%  '
$ws.Range("D17").Value = '0.789'
$ws.Range("E17").Value = '  -3.22%  '
$ws.Range("D18").Value = '42.997.06'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").Value = '13.28'
$ws.Range("E19").Value = '  +6.03%  '
$ws.Range("D20").Value = '6.20'
$ws.Range("E20").Value = '  +0.32%  '
$ws.Range("D21").Value = '0.0₃0908'
$ws.Range("E21").Value = '  +0.21%  '
$ws.Range("D22").Value = '68.08'
$ws.Range("E22").Value = '  +0.47%  '
$ws.Range("D23").Value = '241.17'
$ws.Range("E23").Value = '  +1.64%  '
$ws.Range("E24").Value = '  -2.96%  '
$ws.Range("E25").Value = '  -1.15%  '
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("E27").Value = '  +1.05%  '
$ws.Range("D28").Value = '169.09'
$ws.Range("E28").Value = '  +0.58%  '
$ws.Range("E29").Value = '  -2.68%  '
$ws.Range("D30").Value = '9.17'
$ws.Range("E30").Value = '  -0.16%  '
$ws.Range("D31").Value = '33.44'
$ws.Range("E31").Value = '  -2.29%  '
$ws.Range("D32").Value = '4.96'
$ws.Range("E32").Value = '  +6.43%  '
$ws.Range("D33").Value = '5.19'
$ws.Range("E33").Value = '  +3.17%  '
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").Value = '18.35'
$ws.Range("E35").Value = '  +7.21%  '
$ws.Range("E36").Value = '  -0.75%  '
$ws.Range("D37").Value = '0.0694'
$ws.Range("E37").Value = '  +0.28%  '
$ws.Range("E38").Value = '  +0.58%  '
$ws.Range("D39").Value = '1.81'
$ws.Range("E39").Value = '  +1.03%  '
$ws.Range("D40").Value = '2.76'
$ws.Range("E40").Value = '  -2.41%  '
$ws.Range("E41").Value = '  -0.53%  '
$ws.Range("D42").Value = '1.995.64'
$ws.Range("E42").Value = '  -0.46%  '
$ws.Range("E43").Value = '  +1.05%  '
$ws.Range("E44").Value = '  -0.92%  '
$ws.Range("D45").Value = '17.59'
$ws.Range("E45").Value = '  +0.48%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '2.04'
$ws.Range("E46").Value = '  -12.60%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '2.84'
$ws.Range("E47").Value = '  -0.53%  '
$ws.Range("D48").Value = '76.15'
$ws.Range("E48").Value = '  +8.55%  '
$ws.Range("D49").Value = '54.76'
$ws.Range("E49").Value = '  -1.77%  '
$ws.Range("D50").Value = '2.539.24'
$ws.Range("E50").Value = '  +0.67%  '
$ws.Range("E51").Value = '  +0.13%  '

# Restore the default "Normal" style on the cells we forced to Text so no
# stray number-format style survives on the saved cell (matches original).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
